$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Not stressful"
$ws.Range("B3").Value = "Not stressful"
$ws.Range("B4").Value = "Moderately stressful"
$ws.Range("B5").Value = "Not stressful"
$ws.Range("B6").Value = "Moderately stressful"
$ws.Range("B7").Value = "Moderately stressful"
